# Add team-member attributions to the Gantt task names.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value  = "PROGETTAZIONE (Crepaldi e Malachin)"
$ws.Range("A5").Value  = "SETUP PROGETTO (Malachin)"
$ws.Range("A7").Value  = "SCRITTURA CODICE (Crepaldi)"
$ws.Range("A9").Value  = "GESTIONE ERRORE (Crepaldi)"
$ws.Range("A11").Value = "TEST E DEBUG (Malachin)"
$ws.Range("A14").Value = "CORREZIONE FINALI (Crepaldi)"

# Re-apply the frozen panes (row/col split at E3) and leave the cursor on
# B17, matching the view state the workbook was left in.
$ws.Range("E3").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B17").Select() | Out-Null
